# test code generation module - update evaluations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Execution metrics section updates
# Compilation success -> "no", with a note added
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Wrong method called"

# Runtime without error -> value cleared (no longer "yes")
$ws.Range("B6").ClearContents()

# Assertion validity -> value + note cleared (no longer "yes" / note removed)
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Code BLEU score + detail note updated
$ws.Range("B12").Value = 0.2529545045623279
$ws.Range("C12").Value = "{'codebleu': 0.25295450456232793, 'ngram_match_score': 0.09645012114223925, 'weighted_ngram_match_score': 0.11084301022924441, 'syntax_match_score': 0.5692307692307692, 'dataflow_match_score': 0.23529411764705882}"

# Update the active selection to reflect where the user last clicked
$ws.Activate()
$ws.Range("C7").Select()
